$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 5467
$ws.Range("E3").Value = 17848
$ws.Range("E4").Value = 9798
$ws.Range("E5").Value = 7364
$ws.Range("E6").Value = 16486
$ws.Range("E7").Value = 19686
$ws.Range("E8").Value = 5812
$ws.Range("E9").Value = 13069
$ws.Range("E10").Value = 1434
$ws.Range("E11").Value = 2046
$ws.Range("E12").Value = 3785
$ws.Range("E13").Value = 15853
